# Revert "Increased Slugs and Buckshot damages"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Slugs (ammo_23x75_barrikada) DMG: 2.7 -> 2.5
$ws.Range("H38").Value = 2.5

# Buckshot (ammo_12x70_buck) DMG formula: =9*0.42 -> =9*0.4
$ws.Range("H39").Formula = "=9*0.4"

# restore the selected cell as recorded in the sheet view
$ws.Range("N22").Select()
